$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.593.29"
$ws.Range("E2").Value = "  -5.28%  "

# Row 3
$ws.Range("D3").Value = "3.505.12"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.48%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "390.70"
$ws.Range("E5").Value = "  -6.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.43"
$ws.Range("E6").Value = "  -6.24%  "

# Row 7
$ws.Range("D7").Value = "3.490.89"
$ws.Range("E7").Value = "  -1.96%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("E8").Value = "  -10.39%  "

# Row 9
$ws.Range("E9").Value = "  +0.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.669"
$ws.Range("E10").Value = "  -13.05%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  -14.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000336"
$ws.Range("E12").Value = "  +0.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.50"
$ws.Range("E13").Value = "  -9.02%  "

# Row 14
$ws.Range("D14").Value = "4.055.63"
$ws.Range("E14").Value = "  -1.99%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.14"
$ws.Range("E15").Value = "  -8.43%  "

# Row 16
$ws.Range("E16").Value = "  -3.21%  "

# Row 17
$ws.Range("D17").Value = "3.494.60"
$ws.Range("E17").Value = "  -2.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.59"
$ws.Range("E18").Value = "  -9.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.43"
$ws.Range("E19").Value = "  +0.75%  "

# Row 20
$ws.Range("D20").Value = "63.628.34"
$ws.Range("E20").Value = "  -5.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"
$ws.Range("E21").Value = "  -11.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.83"
$ws.Range("E22").Value = "  -14.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.77"
$ws.Range("E23").Value = "  +2.36%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.18"
$ws.Range("E24").Value = "  -8.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.83"
$ws.Range("E25").Value = "  -8.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.23"
$ws.Range("E26").Value = "  +7.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "32.90"
$ws.Range("E27").Value = "  -5.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.96"
$ws.Range("E28").Value = "  -12.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("E29").Value = "  -13.96%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.90"
$ws.Range("E30").Value = "  -3.73%  "

# Row 31
$ws.Range("E31").Value = "  -6.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.109"
$ws.Range("E32").Value = "  -7.80%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.78"
$ws.Range("E33").Value = "  -8.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.151"
$ws.Range("E34").Value = "  -6.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "36.38"
$ws.Range("E36").Value = "  -11.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.82"
$ws.Range("E37").Value = "  -4.82%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0434"
$ws.Range("E38").Value = "  -12.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.18%  "

# Row 40
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("E40").Value = "  +14.36%  "

# Row 41
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0631"
$ws.Range("E41").Value = "  -13.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.130"
$ws.Range("E42").Value = "  -11.36%  "

# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.03"
$ws.Range("E43").Value = "  +12.82%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "139.16"
$ws.Range("E44").Value = "  -6.64%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -8.19%  "

# Row 46
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.06"
$ws.Range("E46").Value = "  -6.24%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  -12.31%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.00"
$ws.Range("E49").Value = "  -7.40%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.12"
$ws.Range("E50").Value = "  +11.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.275"
$ws.Range("E51").Value = "  -11.80%  "
